$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) cells to Text format so numeric-looking strings
# like "1.000" / "4.250" / "29.235.13" are preserved exactly as typed,
# matching how the source data is stored (plain text, dot as thousands
# separator in some rows).
$textCells = @(
    "D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11",
    "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21",
    "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31",
    "D32", "D33", "D37", "D38", "D40", "D42", "D43", "D44", "D45", "D46",
    "D47", "D48", "D49", "D50", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated price / volume(1h) figures from the latest GitHub Actions refresh run.
$ws.Range("D2").Value = "29.235.13"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "1.862.83"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "0.7134"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("D6").Value = "240.73"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.07725"
$ws.Range("E8").Value = "  -1.18%  "
$ws.Range("D9").Value = "0.3083"
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("D10").Value = "24.95"
$ws.Range("D11").Value = "0.08333"
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("D12").Value = "1.882.01"
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("D13").Value = "5.203"
$ws.Range("E13").Value = "  -1.75%  "
$ws.Range("D14").Value = "0.7139"
$ws.Range("E14").Value = "  -2.11%  "
$ws.Range("D15").Value = "90.98"
$ws.Range("E15").Value = "  -0.47%  "
$ws.Range("D16").Value = "29.249.98"
$ws.Range("E16").Value = "  -0.84%  "
$ws.Range("D17").Value = "5.942"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "242.74"
$ws.Range("E18").Value = "  -2.41%  "
$ws.Range("D19").Value = "0.000007832"
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("D20").Value = "2.129.90"
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("D21").Value = "13.17"
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "7.891"
$ws.Range("E23").Value = "  -1.07%  "
$ws.Range("D24").Value = "0.9994"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "0.1601"
$ws.Range("E25").Value = "  +2.14%  "
$ws.Range("D26").Value = "163.25"
$ws.Range("D27").Value = "8.896"
$ws.Range("E27").Value = "  -1.86%  "
$ws.Range("D28").Value = "18.51"
$ws.Range("E28").Value = "  +0.84%  "
$ws.Range("D29").Value = "1.344"
$ws.Range("E29").Value = "  -1.56%  "
$ws.Range("D30").Value = "1.498"
$ws.Range("E30").Value = "  +0.95%  "
$ws.Range("D31").Value = "4.419"
$ws.Range("E31").Value = "  +0.80%  "
$ws.Range("D32").Value = "4.250"
$ws.Range("E32").Value = "  +2.30%  "
$ws.Range("D33").Value = "0.8392"
$ws.Range("E33").Value = "  +16.02%  "
$ws.Range("E34").Value = "  -2.45%  "
$ws.Range("E35").Value = "  -0.86%  "
$ws.Range("E36").Value = "  -2.82%  "
$ws.Range("D37").Value = "2.678"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").Value = "0.01852"
$ws.Range("E38").Value = "  -0.68%  "
$ws.Range("E39").Value = "  -1.24%  "
$ws.Range("D40").Value = "1.172.97"
$ws.Range("E40").Value = "  -5.05%  "
$ws.Range("E41").Value = "  +1.67%  "
$ws.Range("D42").Value = "0.8972"
$ws.Range("E42").Value = "  -0.87%  "
$ws.Range("D43").Value = "72.84"
$ws.Range("E43").Value = "  -1.29%  "
$ws.Range("D44").Value = "0.9993"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").Value = "102.17"
$ws.Range("E45").Value = "  -1.64%  "
$ws.Range("D46").Value = "2.027.65"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").Value = "0.5184"
$ws.Range("D48").Value = "1.787"
$ws.Range("E48").Value = "  +1.25%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "9.332"
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").Value = "7.055"
$ws.Range("E50").Value = "  -0.46%  "
$ws.Range("B51").Value = "Frax"
$ws.Range("C51").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D51").Value = "1.000"
$ws.Range("E51").Value = "  +0.11%  "
